$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Mdk"
$ws.Range("C2").Value = "Ptprz1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.891504666666667
$ws.Range("H2").Value = 8.674514
$ws.Range("I2").Value = 0.1213590456377548
$ws.Range("J2").Value = 0.1213590456377548
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.03995766666666666
$ws.Range("N2").Value = 0.119873
$ws.Range("O2").Value = 0.005314930928687666
$ws.Range("P2").Value = 0.005314930928687667
$ws.Range("Q2").Value = 0.1155377796357778
$ws.Range("R2").Value = 1.039840016722
$ws.Range("S2").Value = 0.0006450149451361207
$ws.Range("T2").Value = 0.0006450149451361208

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Mdk"
$ws.Range("C3").Value = "Ptprz1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.891504666666667
$ws.Range("H3").Value = 8.674514
$ws.Range("I3").Value = 0.1213590456377548
$ws.Range("J3").Value = 0.1213590456377548
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.03069133333333333
$ws.Range("N3").Value = 0.092074
$ws.Range("O3").Value = 0.004082378436578614
$ws.Range("P3").Value = 0.004082378436578615
$ws.Range("Q3").Value = 0.08874413355955556
$ws.Range("R3").Value = 0.798697202036
$ws.Range("S3").Value = 0.00049543355099533
$ws.Range("T3").Value = 0.0004954335509953301

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Mdk"
$ws.Range("C4").Value = "Ptprz1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.891504666666667
$ws.Range("H4").Value = 8.674514
$ws.Range("I4").Value = 0.1213590456377548
$ws.Range("J4").Value = 0.1213590456377548
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 7.447354000000001
$ws.Range("N4").Value = 22.342062
$ws.Range("O4").Value = 0.9906026906347337
$ws.Range("P4").Value = 0.9906026906347338
$ws.Range("Q4").Value = 21.53405884531867
$ws.Range("R4").Value = 193.806529607868
$ws.Range("S4").Value = 0.1202185971416233
$ws.Range("T4").Value = 0.1202185971416233

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Mdk"
$ws.Range("C5").Value = "Ptprz1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 12.04042966666667
$ws.Range("H5").Value = 36.121289
$ws.Range("I5").Value = 0.505347637947847
$ws.Range("J5").Value = 0.505347637947847
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.03995766666666666
$ws.Range("N5").Value = 0.119873
$ws.Range("O5").Value = 0.005314930928687666
$ws.Range("P5").Value = 0.005314930928687667
$ws.Range("Q5").Value = 0.4811074751441111
$ws.Range("R5").Value = 4.329967276297
$ws.Range("S5").Value = 0.002685887790668269
$ws.Range("T5").Value = 0.002685887790668269

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Mdk"
$ws.Range("C6").Value = "Ptprz1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 12.04042966666667
$ws.Range("H6").Value = 36.121289
$ws.Range("I6").Value = 0.505347637947847
$ws.Range("J6").Value = 0.505347637947847
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.03069133333333333
$ws.Range("N6").Value = 0.092074
$ws.Range("O6").Value = 0.004082378436578614
$ws.Range("P6").Value = 0.004082378436578615
$ws.Range("Q6").Value = 0.3695368403762223
$ws.Range("R6").Value = 3.325831563386001
$ws.Range("S6").Value = 0.002063020300134227
$ws.Range("T6").Value = 0.002063020300134227

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Mdk"
$ws.Range("C7").Value = "Ptprz1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 12.04042966666667
$ws.Range("H7").Value = 36.121289
$ws.Range("I7").Value = 0.505347637947847
$ws.Range("J7").Value = 0.505347637947847
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 7.447354000000001
$ws.Range("N7").Value = 22.342062
$ws.Range("O7").Value = 0.9906026906347337
$ws.Range("P7").Value = 0.9906026906347338
$ws.Range("Q7").Value = 89.66934203976868
$ws.Range("R7").Value = 807.0240783579181
$ws.Range("S7").Value = 0.5005987298570445
$ws.Range("T7").Value = 0.5005987298570446

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Mdk"
$ws.Range("C8").Value = "Ptprz1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 8.894099
$ws.Range("H8").Value = 26.682297
$ws.Range("I8").Value = 0.3732933164143983
$ws.Range("J8").Value = 0.3732933164143982
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.03995766666666666
$ws.Range("N8").Value = 0.119873
$ws.Range("O8").Value = 0.005314930928687666
$ws.Range("P8").Value = 0.005314930928687667
$ws.Range("Q8").Value = 0.3553874431423333
$ws.Range("R8").Value = 3.198486988281
$ws.Range("S8").Value = 0.001984028192883276
$ws.Range("T8").Value = 0.001984028192883277

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Mdk"
$ws.Range("C9").Value = "Ptprz1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 8.894099
$ws.Range("H9").Value = 26.682297
$ws.Range("I9").Value = 0.3732933164143983
$ws.Range("J9").Value = 0.3732933164143982
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.03069133333333333
$ws.Range("N9").Value = 0.092074
$ws.Range("O9").Value = 0.004082378436578614
$ws.Range("P9").Value = 0.004082378436578615
$ws.Range("Q9").Value = 0.2729717571086667
$ws.Range("R9").Value = 2.456745813978
$ws.Range("S9").Value = 0.001523924585449057
$ws.Range("T9").Value = 0.001523924585449057

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Mdk"
$ws.Range("C10").Value = "Ptprz1"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 8.894099
$ws.Range("H10").Value = 26.682297
$ws.Range("I10").Value = 0.3732933164143983
$ws.Range("J10").Value = 0.3732933164143982
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 7.447354000000001
$ws.Range("N10").Value = 22.342062
$ws.Range("O10").Value = 0.9906026906347337
$ws.Range("P10").Value = 0.9906026906347338
$ws.Range("Q10").Value = 66.23750376404601
$ws.Range("R10").Value = 596.1375338764141
$ws.Range("S10").Value = 0.369785363636066
$ws.Range("T10").Value = 0.369785363636066

